# Test Suite for 3 asmts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename the group from Ren_GP_0001 -> Ren_GP_0005
$ws.Range("B2").Value = "Ren_GP_0005"

# Row 3 used to hold a second data row (Ren_GP_0002 / Percentage / System Formed).
# That row is removed, leaving only the date cell (C3) with its formatting but no value.
$ws.Range("A3:G3").ClearContents()

# Update the selected/active cell shown when the workbook is reopened.
$ws.Range("B2").Select()
